$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.084.35"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.902.58"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.73"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5069"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3925"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09271"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.136"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.91"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.367"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.81"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.901.84"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.297"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001120"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.41"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06576"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.213"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.122.72"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.318"
$ws.Range("E25").Value = "  +1.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.608"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.119.27"
$ws.Range("E27").Value = "  +1.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.94"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.49"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.27"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.087"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1069"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.602"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.611"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.590"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06673"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02405"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2181"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.230"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.261"
$ws.Range("E40").Value = "  +7.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6367"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.985"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.44"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5986"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.705"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.277"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.007"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.79"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.178"
$ws.Range("E51").Value = "  -1.35%  "
